$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor cell: untouched, style 0 (General), used to restore formatting
# after writing numeric-looking text (which Excel would otherwise mark
# with a quote-prefix style when entered via apostrophe).
$donor = $ws.Range("B12")

# Map of cell -> new text value (all values are stored as TEXT, matching
# the workbook's original shared-string representation, not as numbers).
$updates = @{
    "B13" = "37.68"
    "C13" = "1.77"
    "D13" = "39.46"
    "B14" = "36.19"
    "C14" = "32.54"
    "D14" = "68.73"
    "B16" = "95.31"
    "C16" = "4.49"
    "B20" = "16.42"
    "C20" = "34.88"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Leading apostrophe forces the value to be stored as text rather than
    # being reinterpreted as a number.
    $cell.Value = "'" + $updates[$addr]
    # Re-apply the original (General) formatting so the cell's style index
    # doesn't drift away from the rest of the table.
    $donor.Copy()
    $cell.PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
